# Update the data-source table on the active sheet: rename a couple of
# column headers, append extra values to a couple of column cells, clear
# a stray leftover date in the "HIS" date column, and move the active
# selection — matching the upstream "use multiprocess for data processing"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells D1 ("C" -> "中文列测试") and G1 ("HIS" -> "时间").
$ws.Range("D1").Value = "中文列测试"
$ws.Range("G1").Value = "时间"

# Extend a couple of existing cell values with extra comma separated items.
$ws.Range("C15").Value = "sdg,dasf,fas"
$ws.Range("C10").Value = "ngfsg,qewr"

# Clear the stray date value left in I14 (style is preserved).
$ws.Range("I14").ClearContents()

# Move the active selection to I24.
$ws.Range("I24").Select()
